$wb = $excel.ActiveWorkbook

# "Greece Market" test data is added the same way the previous market (Croatia) was: by
# duplicating the existing "Croatia" sheet -- which already has the right layout, column
# widths, styles and merged cells -- and dropping it in right after Croatia, at the end of
# the tab strip.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $croatia)
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Fill in the Greece-specific data. (Set B4 before B2 so the new shared strings come out in
# the same order as the rest of the workbook: the processor/repeater code first, then the
# market name.)
$greece.Range("B4").Value = "NGC-4119/T3166"
$greece.Range("B2").Value = "Greece Market"

# Croatia is no longer the focused tab, so its old single-cell selection turns into a
# whole-sheet selection (what Excel records once a tab is deselected after the user had
# clicked a row/column header there).
$croatia.Range("A1:XFD1048576").Select()

# Greece becomes the new active tab, with B10 selected.
$greece.Activate()
$greece.Range("B10").Select()
